# "big fix to smoothness parameter"
# - bump the cmax (smoothness) value in C6
# - widen the tab-bar/scrollbar split (cosmetic window chrome)
# - append four new camparam rows: trimedgeof, openradius, closewidth, closeheight

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen the sheet-tab area relative to the horizontal scrollbar.
$excel.ActiveWindow.TabRatio = 713

# Smoothness parameter fix: cmax goes from 6500 to 7000.
$ws.Range("C6").Value = 7000

# New parameter rows appended after medfiltsize (row 18).
$ws.Range("A19").Value = "trimedgeof"
$ws.Range("B19").Value = 5
$ws.Range("C19").Value = 5

$ws.Range("A20").Value = "openradius"
$ws.Range("B20").Value = 4
$ws.Range("C20").Value = 4

$ws.Range("A21").Value = "closewidth"
$ws.Range("B21").Value = 5
$ws.Range("C21").Value = 5

$ws.Range("A22").Value = "closeheight"
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = 1

# Leave the selection where the author left it.
$ws.Range("C22").Select() | Out-Null
